$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting existing rows 3-12 down to 4-13.
$ws.Rows.Item(3).Insert()

# Copy the formatting of row 2 (A2:I2) into the freshly inserted row 3 (same
# row style as the "Aus bus" / Western Canada summary row immediately above
# it), without touching the rest of the (unbounded) row.
$ws.Range("A2:I2").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122)

# Populate the new row 3 with the Canada entry.
$ws.Cells.Item(3, 1).Value = "Aus bus"
$ws.Cells.Item(3, 4).Value = "Canada"
$ws.Cells.Item(3, 5).Value = "Saskatchewan"
$ws.Cells.Item(3, 7).Value = 11

$ws.Range("E3").Select()
